# Rebuild the button/comment/row table on the active sheet:
#  - Column B ("Comment") values are cleared (all were "N/A", now blank)
#  - New sign types are inserted and the "Row" grouping numbers are
#    reassigned so kph signs = 1, highway/urban signs = 2, other signs = 3
#  - The table grows from 19 data rows to 25 (header + 24 entries)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old table footprint first (old range was A1:C19) so no
# stray cells remain below/beside the new, larger table.
$ws.Range("A1:C19").Clear() | Out-Null

# Seed the brand-new button names first, in the order they were originally
# typed, so new shared-string entries land in that same order.
$ws.Range("A25").Value = "Yield sign"
$ws.Range("A19").Value = "End of highway"
$ws.Range("A16").Value = "Unlimited speed sign"
$ws.Range("A20").Value = "End of urban/city"
$ws.Range("A17").Value = "Start of highway"
$ws.Range("A18").Value = "Start of urban/city"

$data = @(
    @("Button name", "Comment", "Row"),
    @("10 kph sign", $null, 1),
    @("20 kph sign", $null, 1),
    @("30 kph sign", $null, 1),
    @("40 kph sign", $null, 1),
    @("50 kph sign", $null, 1),
    @("60 kph sign", $null, 1),
    @("70 kph sign", $null, 1),
    @("80 kph sign", $null, 1),
    @("90 kph sign", $null, 1),
    @("100 kph sign", $null, 1),
    @("110 kph sign", $null, 1),
    @("120 kph sign", $null, 1),
    @("130 kph sign", $null, 1),
    @("140 kph sign", $null, 1),
    @("Unlimited speed sign", $null, 1),
    @("Start of highway", $null, 2),
    @("Start of urban/city", $null, 2),
    @("End of highway", $null, 2),
    @("End of urban/city", $null, 2),
    @("Stop sign", $null, 3),
    @("Traffic light", $null, 3),
    @("Roundabout", $null, 3),
    @("Speed bump", $null, 3),
    @("Yield sign", $null, 3)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]

    if ($row[1]) {
        $ws.Cells.Item($r, 2).Value = $row[1]
    } else {
        $ws.Cells.Item($r, 2).Value = $null
    }

    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Column A got wider to fit the new, longer sign names (best-fit for the
# longest entry, "Unlimited speed sign").
$ws.Columns.Item(1).ColumnWidth = 18.75

# Move the active selection like the author's session ended up.
$ws.Range("G26").Select() | Out-Null
